# Update the "Output" worksheet's row 2 with the latest flight-selection
# test run results (new timestamp + fare figures).
#
# This mirrors how the test automation tool writes results: it first
# stamps the row with a timestamp, then fills in the fare columns, and
# finally re-stamps the row with the timestamp of when the run actually
# completed (overwriting the earlier timestamp value).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

# Interim timestamp written first (creates a shared-string entry that is
# later superseded once the final timestamp is recorded below).
$ws.Range("A2").Value = "25/01/2022 1:51:28 pm"

# Updated fare figures for the row.
$ws.Range("D2").Value = "₹1,21,120"
$ws.Range("E2").Value = "₹8,328"
$ws.Range("G2").Value = "₹1,29,458"

# Final timestamp once the run completed.
$ws.Range("A2").Value = "25/01/2022 2:23:05 pm"
